# Market_data.xlsx edit — "Added 3.1 - Simulation for complex strategies, with Stop Loss.
# Processed backtest and generated summaries for complex strategies"
#
# Net effect observed in the target OOXML diff:
#  - The "INJ" row and the "NVDA" row are removed from the data table (rows shift up).
#  - Header row translated from Polish to English: "Poczatek"->"Start", "Koniec"->"End",
#    "Plik"->"File" (columns B, C, E respectively; column order unchanged).
#  - Column B (Start date) widened to display width 25.
#  - Selection moved to row 23 (entire row) and the view scrolled down toward row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Remove the two tickers that were dropped from the dataset -------------
# Find "INJ" and "NVDA" in column A (tickers are unique) and delete their rows.
$lookFor = @("INJ", "NVDA")
foreach ($ticker in $lookFor) {
    $found = $ws.Columns("A:A").Find($ticker, [Type]::Missing, [Type]::Missing, 1)
    if ($found -ne $null) {
        $found.EntireRow.Delete()
    }
}

# --- Re-label the header row (Polish -> English) ----------------------------
$ws.Range("B1").Value = "Start"
$ws.Range("C1").Value = "End"
$ws.Range("E1").Value = "File"

# --- Widen column B (Start date) --------------------------------------------
$ws.Columns(2).ColumnWidth = 24.1666666666667

# --- Update the view: select row 23 and scroll so row 25 is at the top -----
$ws.Range("A23:XFD23").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
